$d = $word.ActiveDocument

$replacements = @(
    @("191×3=", "285×7="),
    @("557×8=", "181×8="),
    @("683×4=", "484×5="),
    @("624×8=", "998×8="),
    @("547×3=", "162×9="),
    @("564×9=", "430×3="),
    @("481×6=", "920×2="),
    @("108×7=", "781×9="),
    @("116×9=", "222×8="),
    @("308×9=", "517×2="),
    @("964×4=", "997×5="),
    @("834×9=", "215×2="),
    @("420×3=", "117×5="),
    @("274×7=", "891×8="),
    @("254×7=", "578×8="),
    @("530×3=", "824×3="),
    @("267×6=", "913×8="),
    @("518×5=", "109×7="),
    @("965×4=", "845×3="),
    @("826×7=", "239×5="),
    @("174×8=", "995×2="),
    @("755×8=", "123×9="),
    @("645×4=", "426×5="),
    @("670×6=", "800×6="),
    @("347×9=", "400×2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
